# Generate Report for Handoff
#
# The CI status report regenerated this run: the entry for
# "b1e97847-16d4-4e7d-9d35-fa135838df80.md" is now listed first (still
# "In Translation"), and the entry for
# "709df42e-529f-4d3c-a636-7f3ee41b8d46.md" moved to the second row and
# advanced to "Ready for handoff" with refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")

$ws.Range("A2").Value = "b1e97847-16d4-4e7d-9d35-fa135838df80.md"
$ws.Range("A3").Value = "709df42e-529f-4d3c-a636-7f3ee41b8d46.md"
$ws.Range("B3").Value = "Ready for handoff"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "2016-12-11 18:12:52"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")

$ws.Range("A2").Value = "b1e97847-16d4-4e7d-9d35-fa135838df80.md"
$ws.Range("D2").Value = "b1e97847-16d4-4e7d-9d35-fa135838df80.91344925d237b1e3bd4deb0aba58197ecaa0dd56.zh-cn.xlf"

$ws.Range("A3").Value = "709df42e-529f-4d3c-a636-7f3ee41b8d46.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "709df42e-529f-4d3c-a636-7f3ee41b8d46.d8ce0dc40ca8e472981045caf1beaec00f6b9395.zh-cn.xlf"
$ws.Range("E3").Value = "2016-03-11 18:12:49"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")

$ws.Range("A2").Value = "b1e97847-16d4-4e7d-9d35-fa135838df80.md"
$ws.Range("D2").Value = "b1e97847-16d4-4e7d-9d35-fa135838df80.91344925d237b1e3bd4deb0aba58197ecaa0dd56.de-de.xlf"

$ws.Range("A3").Value = "709df42e-529f-4d3c-a636-7f3ee41b8d46.md"
$ws.Range("C3").Value = "Ready for handoff"
$ws.Range("D3").Value = "709df42e-529f-4d3c-a636-7f3ee41b8d46.d8ce0dc40ca8e472981045caf1beaec00f6b9395.de-de.xlf"
$ws.Range("E3").Value = "2016-03-11 18:12:52"
